$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap prolificid/name/B between row 4 and row 5 (Jennifer/Maggie records trade places)
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = "5f2c1a97a6809c060fec8820"
$ws.Cells.Item(4, 4).Value = "Maggie"

$ws.Cells.Item(5, 2).Value = 10
$ws.Cells.Item(5, 3).Value = "60a71d27a66fac796ad4de6f"
$ws.Cells.Item(5, 4).Value = "Jennifer"

# Update realeffort (F) ranking score values for all data rows
$ws.Cells.Item(2, 6).Value = 11.35146450363736
$ws.Cells.Item(3, 6).Value = 10.46509312749224
$ws.Cells.Item(4, 6).Value = 8.345780257992518
$ws.Cells.Item(5, 6).Value = 8.335815743434177
$ws.Cells.Item(6, 6).Value = 7.149319585641714
$ws.Cells.Item(7, 6).Value = 6.428054772178267
$ws.Cells.Item(8, 6).Value = 6.387792560901622
$ws.Cells.Item(9, 6).Value = 5.358859522459975
$ws.Cells.Item(10, 6).Value = 5.046880051605173
$ws.Cells.Item(11, 6).Value = 4.123973324417245
$ws.Cells.Item(12, 6).Value = 2.400210451344222
$ws.Cells.Item(13, 6).Value = 0.1240641252646651
